# Generate Report for handoff
# Updates the "Latest Handoff Datetime" column (D) for row 4 on the
# zh-cn and de-de sheets, reflecting a newly generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-13 04:35:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-13 04:35:46"
